# create pph 21 manfaat
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows (10-13) ---

# Row 10
$ws.Range("E10").Value = 38
$ws.Range("H10").Value = 1000001
$ws.Range("I10").Value = "ARDI"

# Row 11
$ws.Range("E11").Value = 39
$ws.Range("H11").Value = 1000002
$ws.Range("I11").Value = "DANI"

# Row 12
$ws.Range("E12").Value = 40
$ws.Range("H12").Value = 1000003
$ws.Range("I12").Value = "ALIM"

# Row 13
$ws.Range("D13").Value = "AIG"
$ws.Range("E13").Value = 41
$ws.Range("H13").Value = 1000004
$ws.Range("I13").Value = "ALIM"
$ws.Range("N13").Value = "lagi3"

# --- Insert new row 14 with same formatting as row 13, then populate ---
$ws.Rows.Item(14).Insert()

$ws.Range("B14").Value = 5
$ws.Range("C14").Value = "Extra Manfaat 5"
$ws.Range("D14").Value = "AIK"
$ws.Range("E14").Value = 42
$ws.Range("F14").Value = 40948
$ws.Range("G14").Value = 222
$ws.Range("H14").Value = 1000005
$ws.Range("I14").Value = "ALIM"
$ws.Range("J14").Value = 43957
$ws.Range("K14").Value = 43977
$ws.Range("L14").Value = 222222
$ws.Range("M14").Value = 444444
$ws.Range("N14").Value = "lagi3"

# --- Update sheet view / selection (scrolled to column C, row 2; H14 selected) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("H14").Select()
